$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"="10.06754901453273"; "C"="8.719607370883393"; "E"="11.36490160194159"; "F"="16.86991607391233"; "G"="3.709809526601005"; "I"="34.50677836509271"; "K"="11.66667725932309"; "M"="15.13666493759706" }
  3 = @{ "B"="9.968447506283031"; "C"="8.487629773458554"; "E"="11.11966213943457"; "F"="15.89584955866808"; "G"="3.713001074596066"; "I"="34.01730367309228"; "K"="11.58148529128995"; "M"="15.01758133747354" }
  4 = @{ "B"="9.911542330116175"; "C"="8.344155149006427"; "E"="10.96990660979174"; "F"="15.26997757108491"; "G"="3.715059317752173"; "I"="33.71461300459328"; "K"="11.53274615658138"; "M"="14.948503449272" }
  5 = @{ "B"="9.889375733596411"; "C"="8.285516874005507"; "E"="10.90917839326414"; "F"="15.00819731993403"; "G"="3.71592296312655"; "I"="33.59081601622951"; "K"="11.51380205592976"; "M"="14.9213957231947" }
  6 = @{ "B"="9.885757616459049"; "C"="8.275772405198785"; "E"="10.8991152390415"; "F"="14.96433081551593"; "G"="3.716067877167042"; "I"="33.57023527878219"; "K"="11.51071236602853"; "M"="14.91695816913091" }
  7 = @{ "B"="9.911239203503689"; "C"="8.343364903837012"; "E"="10.96908627780457"; "F"="15.26647399323137"; "G"="3.715070864252758"; "I"="33.71294512653242"; "K"="11.53248692987003"; "M"="14.94813361264761" }
  8 = @{ "B"="10.03257973417685"; "C"="8.639890056792892"; "E"="11.28022681322545"; "F"="16.53996406344768"; "G"="3.710889561652401"; "I"="34.33850949594751"; "K"="11.63657486400214"; "M"="15.09478277425417" }
  9 = @{ "B"="10.30033527830636"; "C"="9.209387111467866"; "E"="11.89296554344761"; "F"="19.00274580682531"; "G"="3.703468160889545"; "I"="35.54413985252182"; "K"="11.86804562903271"; "M"="15.41316874998257" }
  10 = @{ "B"="10.5131000302429"; "C"="9.615830362476437"; "E"="12.33964868823459"; "F"="20.67494806633232"; "G"="3.698483870657891"; "I"="36.41176143162336"; "K"="12.05338264713928"; "M"="15.66407508942884" }
  11 = @{ "B"="10.61289404178912"; "C"="9.797256347496397"; "E"="12.54109368443806"; "F"="21.3917225636224"; "G"="3.696316740015361"; "I"="36.80143817841975"; "K"="12.14068889164146"; "M"="15.78149616883085" }
  12 = @{ "B"="10.65107419823271"; "C"="9.865393617183072"; "E"="12.61704575545672"; "F"="21.65686569030329"; "G"="3.695510418602856"; "I"="36.94819057604216"; "K"="12.17415148443193"; "M"="15.82639636839198" }
  13 = @{ "B"="10.64283473010707"; "C"="9.850745151693383"; "E"="12.60070403014583"; "F"="21.60004134736742"; "G"="3.695683438739993"; "I"="36.91662222696577"; "K"="12.16692732591936"; "M"="15.81670753404353" }
  14 = @{ "B"="10.61602756504844"; "C"="9.80287372063945"; "E"="12.54734936956994"; "F"="21.4136618050453"; "G"="3.696250116918183"; "I"="36.81352831649984"; "K"="12.14343399979206"; "M"="15.78518162190055" }
  15 = @{ "B"="10.59965698324337"; "C"="9.773475662402221"; "E"="12.5146227707788"; "F"="21.29868154950795"; "G"="3.696599086349656"; "I"="36.75027226444022"; "K"="12.12909508398235"; "M"="15.7659266906741" }
  16 = @{ "B"="10.50663517250889"; "C"="9.603898187625946"; "E"="12.32644166941583"; "F"="20.62722412089977"; "G"="3.698627507149665"; "I"="36.38618663499878"; "K"="12.0477347823554"; "M"="15.65646428403374" }
  17 = @{ "B"="10.45030956471284"; "C"="9.498930946809855"; "E"="12.21049210828746"; "F"="20.20408069617459"; "G"="3.699897488237853"; "I"="36.16148445172671"; "K"="11.99856960826214"; "M"="15.59012756032786" }
  18 = @{ "B"="10.41819807892643"; "C"="9.438233824787282"; "E"="12.14364035767374"; "F"="19.95656407809808"; "G"="3.700637389044169"; "I"="36.03177543540056"; "K"="11.97057485059966"; "M"="15.5522836484946" }
  19 = @{ "B"="10.40737599044783"; "C"="9.417629503137649"; "E"="12.12098038956472"; "F"="19.87204792380562"; "G"="3.700889531087513"; "I"="35.98778088506106"; "K"="11.96114591627744"; "M"="15.53952490560793" }
  20 = @{ "B"="10.45627627504609"; "C"="9.510138789254647"; "E"="12.22285234638704"; "F"="20.24955283636157"; "G"="3.699761320116161"; "I"="36.18545329706559"; "K"="12.00377416671957"; "M"="15.59715727109626" }
  21 = @{ "B"="10.62389122195313"; "C"="9.816950549020797"; "E"="12.56303050669891"; "F"="21.46857628470577"; "G"="3.696083281781536"; "I"="36.8438321710227"; "K"="12.1503238952346"; "M"="15.79443002579377" }
  22 = @{ "B"="10.73569278786588"; "C"="10.01415086537934"; "E"="12.78339554880726"; "F"="22.22866616901552"; "G"="3.693762919635391"; "I"="37.26936178864637"; "K"="12.2484295499066"; "M"="15.92587948896553" }
  23 = @{ "B"="10.67582974538388"; "C"="9.909225447822264"; "E"="12.66598673827289"; "F"="21.82633154475857"; "G"="3.694993735222546"; "I"="37.04271284406391"; "K"="12.1958656002232"; "M"="15.85550410389109" }
  24 = @{ "B"="10.45357787841398"; "C"="9.505072808967153"; "E"="12.21726487108054"; "F"="20.22900810905287"; "G"="3.699822851278337"; "I"="36.17461860965795"; "K"="12.00142033962717"; "M"="15.59397822045923" }
  25 = @{ "B"="10.22492249294969"; "C"="9.05710799384266"; "E"="11.727462979975"; "F"="18.34778573295691"; "G"="3.705393180261207"; "I"="35.22079673885643"; "K"="11.80263726289313"; "M"="15.32391772310549" }
}

foreach ($r in $data.Keys) {
  foreach ($c in $data[$r].Keys) {
    $ws.Range("$c$r").Value = [double]$data[$r][$c]
  }
}